# Auto-generated COM-interop script applying the PlayerPerformance_4386 diff
$wb = $excel.ActiveWorkbook

# --- Update 'ODI Batting' sheet: MATCH_CARD_LINK -> MATCH_CODE (D column) ---
$battingWs = $wb.Worksheets.Item("ODI Batting")
$battingWs.Range("D1").Value = "MATCH_CODE"

$battingCodes = @("3898","3923","3924","3927","3929","3931","3937","3940","3942","3945","3947","3950","3966","3967","3968","3972","3973","3975","3977","3981","3984","3988","4032","4035","4041","4067","4069","4071","4074","4076","4108","4115","4123","4125","4166","4167","4168","4169","4170","4222","4224","4226","4564","4565","4567","4597","4600","4601","4660","4663","4666","4725","4728","4732")
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $r = $i + 2
    $battingWs.Range("D$r").Value = "'" + $battingCodes[$i]
}

# Remove the empty placeholder cells in column B for rows 8, 24 and 25
$battingWs.Range("B8").ClearContents()
$battingWs.Range("B24").ClearContents()
$battingWs.Range("B25").ClearContents()

# --- Update 'ODI Bowling' sheet: MATCH_CARD_LINK -> MATCH_CODE (B column) ---
$bowlingWs = $wb.Worksheets.Item("ODI Bowling")
$bowlingWs.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @("3898","3923","3924","3929","3931","3940","3942","3945","3950","3967","3968","3972","3973","3975","3977","3981","3984","4032","4035","4041","4069","4074","4076","4108","4115","4123","4564","4565","4567","4600")
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $r = $i + 2
    $bowlingWs.Range("B$r").Value = "'" + $bowlingCodes[$i]
}

# --- Add new 'Player Info' worksheet (becomes the first sheet) ---
$playerInfoWs = $wb.Worksheets.Add()
$playerInfoWs.Name = "Player Info"

$playerInfoData = New-Object 'object[,]' 2,4
$playerInfoData[0,0] = "ID"
$playerInfoData[0,1] = "NAME"
$playerInfoData[0,2] = "BATTING_HAND"
$playerInfoData[0,3] = "BOWL_STYLE"
$playerInfoData[1,0] = "'4386"
$playerInfoData[1,1] = "Travis Michael Head"
$playerInfoData[1,2] = "Left Handed"
$playerInfoData[1,3] = "Right Arm Off Break"
$playerInfoWs.Range("A1:D2").Value = $playerInfoData
$playerInfoWs.Range("A1:D1").Font.Bold = $true

# --- Add new 'ODI Batting Extra' worksheet (becomes the last sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtraWs = $wb.Worksheets.Add($null, $lastSheet)
$battingExtraWs.Name = "ODI Batting Extra"

$battingExtraData = New-Object 'object[,]' 21,6
$battingExtraData[0,0] = "MATCH_CODE"
$battingExtraData[0,1] = "BATTING_POSITION"
$battingExtraData[0,2] = "NUM_4"
$battingExtraData[0,3] = "NUM_6"
$battingExtraData[0,4] = "PERCENT_RUNS_OF_TOTAL"
$battingExtraData[0,5] = "MAN_OF_MATCH"
$battingExtraData[1,0] = "'4166"
$battingExtraData[1,1] = 2
$battingExtraData[1,2] = "'0"
$battingExtraData[1,3] = "'0"
$battingExtraData[1,4] = "'2.34%"
$battingExtraData[1,5] = "NO"
$battingExtraData[2,0] = "'4167"
$battingExtraData[2,5] = "NO"
$battingExtraData[3,0] = "'4168"
$battingExtraData[3,1] = 2
$battingExtraData[3,2] = "'7"
$battingExtraData[3,3] = "'0"
$battingExtraData[3,4] = "'21.34%"
$battingExtraData[3,5] = "NO"
$battingExtraData[4,0] = "'4169"
$battingExtraData[4,5] = "NO"
$battingExtraData[5,0] = "'4170"
$battingExtraData[5,1] = 2
$battingExtraData[5,2] = "'9"
$battingExtraData[5,3] = "'0"
$battingExtraData[5,4] = "'27.32%"
$battingExtraData[5,5] = "NO"
$battingExtraData[6,0] = "'4222"
$battingExtraData[6,5] = "NO"
$battingExtraData[7,0] = "'4224"
$battingExtraData[7,1] = 2
$battingExtraData[7,2] = "'2"
$battingExtraData[7,3] = "'0"
$battingExtraData[7,4] = "'3.46%"
$battingExtraData[7,5] = "NO"
$battingExtraData[8,0] = "'4226"
$battingExtraData[8,1] = 4
$battingExtraData[8,2] = "'1"
$battingExtraData[8,3] = "'0"
$battingExtraData[8,4] = "'2.14%"
$battingExtraData[8,5] = "NO"
$battingExtraData[9,0] = "'4564"
$battingExtraData[9,1] = 1
$battingExtraData[9,2] = "'2"
$battingExtraData[9,3] = "'3"
$battingExtraData[9,4] = "'32.27%"
$battingExtraData[9,5] = "YES"
$battingExtraData[10,0] = "'4565"
$battingExtraData[10,5] = "NO"
$battingExtraData[11,0] = "'4567"
$battingExtraData[11,1] = 1
$battingExtraData[11,2] = "'0"
$battingExtraData[11,3] = "'0"
$battingExtraData[11,5] = "NO"
$battingExtraData[12,0] = "'4597"
$battingExtraData[12,1] = 4
$battingExtraData[12,2] = "'2"
$battingExtraData[12,3] = "'0"
$battingExtraData[12,4] = "'12.17%"
$battingExtraData[12,5] = "NO"
$battingExtraData[13,0] = "'4600"
$battingExtraData[13,5] = "NO"
$battingExtraData[14,0] = "'4601"
$battingExtraData[14,1] = 6
$battingExtraData[14,2] = "'3"
$battingExtraData[14,3] = "'0"
$battingExtraData[14,4] = "'10.63%"
$battingExtraData[14,5] = "NO"
$battingExtraData[15,0] = "'4660"
$battingExtraData[15,5] = "NO"
$battingExtraData[16,0] = "'4663"
$battingExtraData[16,5] = "NO"
$battingExtraData[17,0] = "'4666"
$battingExtraData[17,5] = "NO"
$battingExtraData[18,0] = "'4725"
$battingExtraData[18,1] = 1
$battingExtraData[18,2] = "'1"
$battingExtraData[18,3] = "'0"
$battingExtraData[18,4] = "'2.66%"
$battingExtraData[18,5] = "NO"
$battingExtraData[19,0] = "'4728"
$battingExtraData[19,1] = 1
$battingExtraData[19,2] = "'10"
$battingExtraData[19,3] = "'0"
$battingExtraData[19,4] = "'42.15%"
$battingExtraData[19,5] = "NO"
$battingExtraData[20,0] = "'4732"
$battingExtraData[20,1] = 1
$battingExtraData[20,2] = "'4"
$battingExtraData[20,3] = "'2"
$battingExtraData[20,4] = "'12.27%"
$battingExtraData[20,5] = "NO"
$battingExtraWs.Range("A1:F21").Value = $battingExtraData
$battingExtraWs.Range("A1:F1").Font.Bold = $true

Write-Output "Edit complete"
